$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.542102
$ws.Range("H2").Value = 16.626306
$ws.Range("I2").Value = 0.2361826998234217
$ws.Range("J2").Value = 0.2361826998234217
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.616972666666666
$ws.Range("N2").Value = 16.850918
$ws.Range("O2").Value = 0.2248624183853938
$ws.Range("P2").Value = 0.2248624183853938
$ws.Range("Q2").Value = 31.12983544987867
$ws.Range("R2").Value = 280.168519048908
$ws.Range("S2").Value = 0.05310861306308613
$ws.Range("T2").Value = 0.05310861306308611

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.542102
$ws.Range("H3").Value = 16.626306
$ws.Range("I3").Value = 0.2361826998234217
$ws.Range("J3").Value = 0.2361826998234217
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 12.20369633333333
$ws.Range("N3").Value = 36.611089
$ws.Range("O3").Value = 0.4885465594374674
$ws.Range("P3").Value = 0.4885465594374673
$ws.Range("Q3").Value = 67.63412985635934
$ws.Range("R3").Value = 608.707168707234
$ws.Range("S3").Value = 0.1153862453973848
$ws.Range("T3").Value = 0.1153862453973848

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.542102
$ws.Range("H4").Value = 16.626306
$ws.Range("I4").Value = 0.2361826998234217
$ws.Range("J4").Value = 0.2361826998234217
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.734464333333333
$ws.Range("N4").Value = 5.203393
$ws.Range("O4").Value = 0.06943523989551367
$ws.Range("P4").Value = 0.06943523989551366
$ws.Range("Q4").Value = 9.612578250695334
$ws.Range("R4").Value = 86.51320425625801
$ws.Range("S4").Value = 0.01639940242140938
$ws.Range("T4").Value = 0.01639940242140938

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.542102
$ws.Range("H5").Value = 16.626306
$ws.Range("I5").Value = 0.2361826998234217
$ws.Range("J5").Value = 0.2361826998234217
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.424464
$ws.Range("N5").Value = 16.273392
$ws.Range("O5").Value = 0.2171557822816253
$ws.Range("P5").Value = 0.2171557822816252
$ws.Range("Q5").Value = 30.062932783328
$ws.Range("R5").Value = 270.566395049952
$ws.Range("S5").Value = 0.05128843894154143
$ws.Range("T5").Value = 0.05128843894154141

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.525638333333333
$ws.Range("H6").Value = 10.576915
$ws.Range("I6").Value = 0.1502489091986426
$ws.Range("J6").Value = 0.1502489091986426
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.616972666666666
$ws.Range("N6").Value = 16.850918
$ws.Range("O6").Value = 0.2248624183853938
$ws.Range("P6").Value = 0.2248624183853938
$ws.Range("Q6").Value = 19.80341415088555
$ws.Range("R6").Value = 178.23072735797
$ws.Range("S6").Value = 0.03378533308217421
$ws.Range("T6").Value = 0.0337853330821742

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3.525638333333333
$ws.Range("H7").Value = 10.576915
$ws.Range("I7").Value = 0.1502489091986426
$ws.Range("J7").Value = 0.1502489091986426
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 12.20369633333333
$ws.Range("N7").Value = 36.611089
$ws.Range("O7").Value = 0.4885465594374674
$ws.Range("P7").Value = 0.4885465594374673
$ws.Range("Q7").Value = 43.02581960115944
$ws.Range("R7").Value = 387.232376410435
$ws.Range("S7").Value = 0.07340358764822928
$ws.Range("T7").Value = 0.07340358764822925

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.525638333333333
$ws.Range("H8").Value = 10.576915
$ws.Range("I8").Value = 0.1502489091986426
$ws.Range("J8").Value = 0.1502489091986426
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.734464333333333
$ws.Range("N8").Value = 5.203393
$ws.Range("O8").Value = 0.06943523989551367
$ws.Range("P8").Value = 0.06943523989551366
$ws.Range("Q8").Value = 6.115093941399445
$ws.Range("R8").Value = 55.035845472595
$ws.Range("S8").Value = 0.010432569054247
$ws.Range("T8").Value = 0.01043256905424699

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.525638333333333
$ws.Range("H9").Value = 10.576915
$ws.Range("I9").Value = 0.1502489091986426
$ws.Range("J9").Value = 0.1502489091986426
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.424464
$ws.Range("N9").Value = 16.273392
$ws.Range("O9").Value = 0.2171557822816253
$ws.Range("P9").Value = 0.2171557822816252
$ws.Range("Q9").Value = 19.12469821618667
$ws.Range("R9").Value = 172.12228394568
$ws.Range("S9").Value = 0.03262741941399212
$ws.Range("T9").Value = 0.0326274194139921

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 13.07613666666666
$ws.Range("H10").Value = 39.22841
$ws.Range("I10").Value = 0.5572537750466107
$ws.Range("J10").Value = 0.5572537750466107
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.616972666666666
$ws.Range("N10").Value = 16.850918
$ws.Range("O10").Value = 0.2248624183853938
$ws.Range("P10").Value = 0.2248624183853938
$ws.Range("Q10").Value = 73.44830224226443
$ws.Range("R10").Value = 661.03472018038
$ws.Range("S10").Value = 0.1253054315113711
$ws.Range("T10").Value = 0.1253054315113711

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 13.07613666666666
$ws.Range("H11").Value = 39.22841
$ws.Range("I11").Value = 0.5572537750466107
$ws.Range("J11").Value = 0.5572537750466107
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 12.20369633333333
$ws.Range("N11").Value = 36.611089
$ws.Range("O11").Value = 0.4885465594374674
$ws.Range("P11").Value = 0.4885465594374673
$ws.Range("Q11").Value = 159.5772010931655
$ws.Range("R11").Value = 1436.19480983849
$ws.Range("S11").Value = 0.272244414532562
$ws.Range("T11").Value = 0.272244414532562

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 13.07613666666666
$ws.Range("H12").Value = 39.22841
$ws.Range("I12").Value = 0.5572537750466107
$ws.Range("J12").Value = 0.5572537750466107
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 1.734464333333333
$ws.Range("N12").Value = 5.203393
$ws.Range("O12").Value = 0.06943523989551367
$ws.Range("P12").Value = 0.06943523989551366
$ws.Range("Q12").Value = 22.68009266612555
$ws.Range("R12").Value = 204.12083399513
$ws.Range("S12").Value = 0.03869304955304202
$ws.Range("T12").Value = 0.03869304955304202

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 13.07613666666666
$ws.Range("H13").Value = 39.22841
$ws.Range("I13").Value = 0.5572537750466107
$ws.Range("J13").Value = 0.5572537750466107
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 5.424464
$ws.Range("N13").Value = 16.273392
$ws.Range("O13").Value = 0.2171557822816253
$ws.Range("P13").Value = 0.2171557822816252
$ws.Range("Q13").Value = 70.93103260741333
$ws.Range("R13").Value = 638.37929346672
$ws.Range("S13").Value = 0.1210108794496356
$ws.Range("T13").Value = 0.1210108794496356

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 1.321440333333333
$ws.Range("H14").Value = 3.964321
$ws.Range("I14").Value = 0.05631461593132515
$ws.Range("J14").Value = 0.05631461593132514
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 5.616972666666666
$ws.Range("N14").Value = 16.850918
$ws.Range("O14").Value = 0.2248624183853938
$ws.Range("P14").Value = 0.2248624183853938
$ws.Range("Q14").Value = 7.422494232964221
$ws.Range("R14").Value = 66.80244809667799
$ws.Range("S14").Value = 0.0126630407287624
$ws.Range("T14").Value = 0.01266304072876239

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 1.321440333333333
$ws.Range("H15").Value = 3.964321
$ws.Range("I15").Value = 0.05631461593132515
$ws.Range("J15").Value = 0.05631461593132514
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 12.20369633333333
$ws.Range("N15").Value = 36.611089
$ws.Range("O15").Value = 0.4885465594374674
$ws.Range("P15").Value = 0.4885465594374673
$ws.Range("Q15").Value = 16.12645655061878
$ws.Range("R15").Value = 145.138108955569
$ws.Range("S15").Value = 0.02751231185929129
$ws.Range("T15").Value = 0.02751231185929128

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 1.321440333333333
$ws.Range("H16").Value = 3.964321
$ws.Range("I16").Value = 0.05631461593132515
$ws.Range("J16").Value = 0.05631461593132514
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.734464333333333
$ws.Range("N16").Value = 5.203393
$ws.Range("O16").Value = 0.06943523989551367
$ws.Range("P16").Value = 0.06943523989551366
$ws.Range("Q16").Value = 2.291991126794777
$ws.Range("R16").Value = 20.627920141153
$ws.Range("S16").Value = 0.003910218866815278
$ws.Range("T16").Value = 0.003910218866815277

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 1.321440333333333
$ws.Range("H17").Value = 3.964321
$ws.Range("I17").Value = 0.05631461593132515
$ws.Range("J17").Value = 0.05631461593132514
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 5.424464
$ws.Range("N17").Value = 16.273392
$ws.Range("O17").Value = 0.2171557822816253
$ws.Range("P17").Value = 0.2171557822816252
$ws.Range("Q17").Value = 7.168105516314666
$ws.Range("R17").Value = 64.512949646832
$ws.Range("S17").Value = 0.01222904447645619
$ws.Range("T17").Value = 0.01222904447645619

